# Atualização contrato 070/25 e hedge
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Update "futuros" curve (drives all the downstream hedge/summary recalcs
#    on Sheet2 (2), hedge, etc. via formulas referencing futuros!B3 etc.)
# ---------------------------------------------------------------------------
$futuros = $wb.Worksheets.Item("futuros")
$futuros.Range("D1").Value2 = 45887

$futuros.Range("B2").Value2  = 340
$futuros.Range("B3").Value2  = 332.4
$futuros.Range("B4").Value2  = 321.25
$futuros.Range("B5").Value2  = 312.95
$futuros.Range("B6").Value2  = 305.45
$futuros.Range("B7").Value2  = 297.4
$futuros.Range("B8").Value2  = 292.55
$futuros.Range("B9").Value2  = 289.5
$futuros.Range("B10").Value2 = 286.35
$futuros.Range("B11").Value2 = 282.9
$futuros.Range("B12").Value2 = 279.1
$futuros.Range("B13").Value2 = 275.75

# ---------------------------------------------------------------------------
# 2) Register new contract 070/25 on "Sheet2": insert a row at 66 (format
#    inherited from row above, matching Excel's default insert behaviour)
#    and fill in the new sale record.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Rows.Item(66).Insert()

$ws2.Range("A66").Value2 = 2025
$ws2.Range("B66").Value2 = "070/25"
$ws2.Range("C66").Value2 = "Mercado Interno"
$ws2.Range("D66").Value2 = "Mundo Café"
$ws2.Range("E66").Value2 = 10
$ws2.Range("F66").Value2 = "Moka"
$ws2.Range("G66").Value2 = "Petrus"
$ws2.Range("H66").Value2 = 0
$ws2.Range("I66").Value2 = 45887
$ws2.Range("J66").Value2 = 4
$ws2.Range("K66").Value2 = 45880
$ws2.Range("P66").Value2 = 2250
$ws2.Range("Q66").Formula = "=E66*P66"

# --- Row 67 (previously row 66): switch its market and re-derive the volume
#     still outstanding for contract 331 from the new split across rows.
$ws2.Range("C67").Value2 = "Mercado Interno"
$ws2.Range("E67").Formula = "=3000-E66-E65"
$ws2.Range("L67").ClearContents()
$ws2.Range("N67").ClearContents()
$ws2.Range("P67").Value2 = 1600
$ws2.Range("Q67").Formula = "=E67*P67"

# --- Row 68 (previously row 67): recompute outstanding volume, update
#     installments and locked-in price.
$ws2.Range("E68").Formula = "=7000-4800-E69-E70"
$ws2.Range("J68").Value2 = 4
$ws2.Range("L68").Value2 = 325

# --- Row 69 (previously row 68): swap in the smaller, earlier-dated lot.
$ws2.Range("E69").Value2 = 64.5
$ws2.Range("I69").Value2 = 45992
$ws2.Range("K69").Value2 = 45992
$ws2.Range("L69").Value2 = 330

# Row 70 (previously row 69) keeps its original values untouched by the
# Insert shift, so nothing further to change there.

# ---------------------------------------------------------------------------
# 3) Extend the autofilter / sort range and the _FilterDatabase defined name
#    to cover the new last row (70).
# ---------------------------------------------------------------------------
$ws2.AutoFilterMode = $false
$ws2.Range("A1:Q70").AutoFilter() | Out-Null

$wb.Names.Item("_xlnm._FilterDatabase").RefersTo = "=Sheet2!`$A`$1:`$Q`$70"

# ---------------------------------------------------------------------------
# 4) Restore view state tweaks captured in the saved workbook.
# ---------------------------------------------------------------------------
$ws2.Activate()
$excel.ActiveWindow.ScrollRow = 5
$ws2.Range("E66").Select() | Out-Null

$hedge = $wb.Worksheets.Item("hedge")
$hedge.Activate()
$hedge.Range("J29").Select() | Out-Null
